# Add a new acronym-key row for "AVLRaPTC" (Annual Vehicle Licensing
# Registration and Property Tax Costs) to the "Key to Variables" sheet,
# directly above the existing "AVMC" row, matching alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row at 190 (pushes "AVMC" and everything below down by one).
# EntireRow.Insert() carries down the formatting (wrap text / fills) from
# the row above, same as Excel does interactively.
$ws.Rows.Item(190).EntireRow.Insert()

$ws.Cells.Item(190, 1).Value = "trans"
$ws.Cells.Item(190, 2).Value = "AVLRaPTC"
$ws.Cells.Item(190, 3).Value = "Annual Vehicle Licensing Registration and Property Tax Costs"
$ws.Cells.Item(190, 6).Value = "medium"

# The long "Meaning" text wraps to two lines at this column width, so the
# row renders taller (30pt) than the single-line rows around it.
$ws.Rows.Item(190).RowHeight = 30
